$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update HUD test data values (L and V columns, rows 2-8)
$ws.Range("L2").Value = 85000
$ws.Range("V2").Value = 30

$ws.Range("L3").Value = 85000
$ws.Range("V3").Value = 100

$ws.Range("V4").Value = 40

$ws.Range("V5").Value = 67

$ws.Range("V6").Value = 45

$ws.Range("V7").Value = 24

$ws.Range("V8").Value = 76

# Update the view: scroll so column I is the top-left visible column,
# and move the active selection to V8
$ws.Application.ActiveWindow.ScrollColumn = 9
$ws.Range("V8").Select()
